$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 451. This shifts all existing rows
# 451-509 down to 452-510, matching the rest of the diff automatically.
$ws.Rows.Item(451).Insert()

# Populate the newly inserted row 451 with its data.
$ws.Cells.Item(451, 1).Value = 8
$ws.Cells.Item(451, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(451, 3).Value = "Coquimbo"
$ws.Cells.Item(451, 4).Value = 45142
$ws.Cells.Item(451, 5).Value = 4
$ws.Cells.Item(451, 6).Value = 100112021
$ws.Cells.Item(451, 7).Value = "Ají"
$ws.Cells.Item(451, 8).Value = "Inferno"
$ws.Cells.Item(451, 9).Value = "Primera"
$ws.Cells.Item(451, 10).Value = 460
$ws.Cells.Item(451, 11).Value = 15000
$ws.Cells.Item(451, 12).Value = 16000
$ws.Cells.Item(451, 13).Value = 15500
$ws.Cells.Item(451, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(451, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(451, 16).Value = 1550
$ws.Cells.Item(451, 17).Value = 10
$ws.Cells.Item(451, 18).Value = "Hortaliza"
